# Auto-generated PowerShell Excel COM-interop script
# Applies numeric cell updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1825
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19299.277
$ws.Range("I32").Value = 13816
$ws.Range("J32").Value = 22040.916
$ws.Range("K32").Value = 13816
$ws.Range("L32").Value = 22040.916
$ws.Range("M32").Value = -13529
$ws.Range("N32").Value = -22614.916
$ws.Range("H61").Value = 1719.4108
$ws.Range("I61").Value = 1123.2766
$ws.Range("J61").Value = 4832.5557
$ws.Range("K61").Value = 1123.2766
$ws.Range("L61").Value = 4832.5557
$ws.Range("M61").Value = -911.2765999999999
$ws.Range("N61").Value = -5256.5557
$ws.Range("H110").Value = 34485090
$ws.Range("I110").Value = 50002230
$ws.Range("J110").Value = 2545.6667
$ws.Range("K110").Value = 50002230
$ws.Range("L110").Value = 2545.6667
$ws.Range("M110").Value = -50000185
$ws.Range("N110").Value = -6635.6667
$ws.Range("H122").Value = 2361.8
$ws.Range("I122").Value = 1808.5
$ws.Range("J122").Value = 4575
$ws.Range("K122").Value = 5425.5
$ws.Range("L122").Value = 13725
$ws.Range("M122").Value = -2975.5
$ws.Range("N122").Value = -18625
$ws.Range("H132").Value = 1374.0741
$ws.Range("I132").Value = 1191.902
$ws.Range("J132").Value = 4471
$ws.Range("K132").Value = 3575.706
$ws.Range("L132").Value = 13413
$ws.Range("M132").Value = -1045.706
$ws.Range("N132").Value = -18473
$ws.Range("H134").Value = 53211
$ws.Range("I134").Value = 53211
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 53211
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -48141
$ws.Range("H136").Value = 1719.4108
$ws.Range("I136").Value = 1123.2766
$ws.Range("J136").Value = 4832.5557
$ws.Range("K136").Value = 3369.8298
$ws.Range("L136").Value = 14497.6671
$ws.Range("M136").Value = -819.8297999999995
$ws.Range("N136").Value = -19597.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 11700
$ws.Range("I4").Value = 4000
$ws.Range("J4").Value = 13625
$ws.Range("K4").Value = 4000
$ws.Range("L4").Value = 13625
$ws.Range("M4").Value = -3885
$ws.Range("N4").Value = -13855
$ws.Range("H105").Value = 1453.7675
$ws.Range("I105").Value = 1142.4
$ws.Range("J105").Value = 2172.3076
$ws.Range("K105").Value = 1142.4
$ws.Range("L105").Value = 2172.3076
$ws.Range("M105").Value = 604.5999999999999
$ws.Range("N105").Value = -5666.3076
$ws.Range("H107").Value = 83668620
$ws.Range("I107").Value = 500430.25
$ws.Range("J107").Value = 250005010
$ws.Range("K107").Value = 500430.25
$ws.Range("L107").Value = 250005010
$ws.Range("M107").Value = -498510.25
$ws.Range("N107").Value = -250008850
$ws.Range("H133").Value = 49275
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49275
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49275
$ws.Range("N133").Value = -59395
$ws.Range("H134").Value = 1648.1111
$ws.Range("I134").Value = 1571.9318
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 4715.7954
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -2180.7954
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 859.1667
$ws.Range("I16").Value = 831.4
$ws.Range("J16").Value = 998
$ws.Range("K16").Value = 831.4
$ws.Range("L16").Value = 998
$ws.Range("M16").Value = -544.4
$ws.Range("N16").Value = -1572
$ws.Range("H22").Value = 467.25
$ws.Range("I22").Value = 399.8
$ws.Range("J22").Value = 579.6667
$ws.Range("K22").Value = 399.8
$ws.Range("L22").Value = 579.6667
$ws.Range("M22").Value = -49.80000000000001
$ws.Range("N22").Value = -1279.6667
$ws.Range("H31").Value = 7841.25
$ws.Range("I31").Value = 10810.934
$ws.Range("J31").Value = 2891.7778
$ws.Range("K31").Value = 10810.934
$ws.Range("L31").Value = 2891.7778
$ws.Range("M31").Value = -10515.934
$ws.Range("N31").Value = -3481.7778
$ws.Range("H34").Value = 7841.25
$ws.Range("I34").Value = 10810.934
$ws.Range("J34").Value = 2891.7778
$ws.Range("K34").Value = 10810.934
$ws.Range("L34").Value = 2891.7778
$ws.Range("M34").Value = -10608.934
$ws.Range("N34").Value = -3295.7778
$ws.Range("H113").Value = 859.1667
$ws.Range("I113").Value = 831.4
$ws.Range("J113").Value = 998
$ws.Range("K113").Value = 831.4
$ws.Range("L113").Value = 998
$ws.Range("M113").Value = 1338.6
$ws.Range("N113").Value = -5338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1487.1111
$ws.Range("I68").Value = 1055.5714
$ws.Range("J68").Value = 2997.5
$ws.Range("K68").Value = 3166.7142
$ws.Range("L68").Value = 8992.5
$ws.Range("M68").Value = -2355.7142
$ws.Range("N68").Value = -10614.5
$ws.Range("H71").Value = 1487.1111
$ws.Range("I71").Value = 1055.5714
$ws.Range("J71").Value = 2997.5
$ws.Range("K71").Value = 9500.142600000001
$ws.Range("L71").Value = 26977.5
$ws.Range("M71").Value = -5444.142600000001
$ws.Range("N71").Value = -35089.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4000073.2
$ws.Range("I3").Value = 3333431
$ws.Range("J3").Value = 6000000
$ws.Range("K3").Value = 3333431
$ws.Range("L3").Value = 6000000
$ws.Range("M3").Value = -3333315
$ws.Range("N3").Value = -6000232
$ws.Range("H6").Value = 643
$ws.Range("I6").Value = 383
$ws.Range("J6").Value = 903
$ws.Range("K6").Value = 383
$ws.Range("L6").Value = 903
$ws.Range("M6").Value = -270
$ws.Range("N6").Value = -1129
$ws.Range("H7").Value = 4836666.5
$ws.Range("I7").Value = 4500000
$ws.Range("J7").Value = 5005000
$ws.Range("K7").Value = 4500000
$ws.Range("L7").Value = 5005000
$ws.Range("M7").Value = -4499888
$ws.Range("N7").Value = -5005224
$ws.Range("H8").Value = 4836666.5
$ws.Range("I8").Value = 4500000
$ws.Range("J8").Value = 5005000
$ws.Range("K8").Value = 4500000
$ws.Range("L8").Value = 5005000
$ws.Range("M8").Value = -4499861
$ws.Range("N8").Value = -5005278
$ws.Range("H10").Value = 1668500
$ws.Range("I10").Value = 1500
$ws.Range("J10").Value = 2502000
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 2502000
$ws.Range("M10").Value = -1331
$ws.Range("N10").Value = -2502338
$ws.Range("H12").Value = 599.6667
$ws.Range("I12").Value = 399.5
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 399.5
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -259.5
$ws.Range("N12").Value = -1280
$ws.Range("H13").Value = 234.5
$ws.Range("I13").Value = 497.5
$ws.Range("J13").Value = 146.83333
$ws.Range("K13").Value = 497.5
$ws.Range("L13").Value = 146.83333
$ws.Range("M13").Value = -358.5
$ws.Range("N13").Value = -424.83333
$ws.Range("H14").Value = 2000494.4
$ws.Range("I14").Value = 3333661.2
$ws.Range("J14").Value = 744
$ws.Range("K14").Value = 3333661.2
$ws.Range("L14").Value = 744
$ws.Range("M14").Value = -3333493.2
$ws.Range("N14").Value = -1080
$ws.Range("H16").Value = 643
$ws.Range("I16").Value = 383
$ws.Range("J16").Value = 903
$ws.Range("K16").Value = 383
$ws.Range("L16").Value = 903
$ws.Range("M16").Value = -133
$ws.Range("N16").Value = -1403
$ws.Range("H17").Value = 333.4
$ws.Range("I17").Value = 286
$ws.Range("J17").Value = 404.5
$ws.Range("K17").Value = 286
$ws.Range("L17").Value = 404.5
$ws.Range("M17").Value = -118
$ws.Range("N17").Value = -740.5
$ws.Range("H33").Value = 15000
$ws.Range("I33").Value = 15000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -14748
$ws.Range("H53").Value = 50000
$ws.Range("I53").Value = 50000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 50000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -49369
$ws.Range("H113").Value = 4124.6665
$ws.Range("I113").Value = 3312
$ws.Range("J113").Value = 5750
$ws.Range("K113").Value = 3312
$ws.Range("L113").Value = 5750
$ws.Range("M113").Value = -1142
$ws.Range("N113").Value = -10090
$ws.Range("H137").Value = 56959
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 56959
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 56959
$ws.Range("N137").Value = -67159

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18521436
$ws.Range("I7").Value = 31252620
$ws.Range("J7").Value = 3348.6365
$ws.Range("K7").Value = 31252620
$ws.Range("L7").Value = 3348.6365
$ws.Range("M7").Value = -31252508
$ws.Range("N7").Value = -3572.6365
$ws.Range("H22").Value = 858.5714
$ws.Range("I22").Value = 863.3333
$ws.Range("J22").Value = 855
$ws.Range("K22").Value = 863.3333
$ws.Range("L22").Value = 855
$ws.Range("M22").Value = -568.3333
$ws.Range("N22").Value = -1445
$ws.Range("H27").Value = 858.5714
$ws.Range("I27").Value = 863.3333
$ws.Range("J27").Value = 855
$ws.Range("K27").Value = 863.3333
$ws.Range("L27").Value = 855
$ws.Range("M27").Value = -756.3333
$ws.Range("N27").Value = -1069
$ws.Range("H55").Value = 437.76
$ws.Range("I55").Value = 288
$ws.Range("J55").Value = 704
$ws.Range("K55").Value = 288
$ws.Range("L55").Value = 704
$ws.Range("M55").Value = -115
$ws.Range("N55").Value = -1050
$ws.Range("H126").Value = 18521436
$ws.Range("I126").Value = 31252620
$ws.Range("J126").Value = 3348.6365
$ws.Range("K126").Value = 93757860
$ws.Range("L126").Value = 10045.9095
$ws.Range("M126").Value = -93755390
$ws.Range("N126").Value = -14985.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3426
$ws.Range("I17").Value = 3426
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 3426
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -3254
$ws.Range("H107").Value = 1388.25
$ws.Range("I107").Value = 1240.6
$ws.Range("J107").Value = 1634.3334
$ws.Range("K107").Value = 3721.8
$ws.Range("L107").Value = 4903.0002
$ws.Range("M107").Value = -1801.8
$ws.Range("N107").Value = -8743.0002
$ws.Range("H126").Value = 2331.9375
$ws.Range("I126").Value = 1966
$ws.Range("J126").Value = 2551.5
$ws.Range("K126").Value = 5898
$ws.Range("L126").Value = 7654.5
$ws.Range("M126").Value = -3428
$ws.Range("N126").Value = -12594.5
